# Updates cryptos list figures (prices / 1h volume %) per the Tue Jul 16 2024
# GitHub Actions refresh commit. Rows 36-39 also swap coin identity (Monero/
# ImmutableX and Maker/Stacks traded ranking positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells (column D) hold values that look like plain numbers
# (e.g. "1.00", "160.78") but must stay literal text, matching the other
# Price cells in the sheet (which are all stored as text, e.g. "65.205.90").
# Force text format first so Excel does not auto-convert them to numbers
# (and drop formatting such as trailing zeros) when .Value is assigned.
$textPriceCells = @(
    "D5", "D6", "D7", "D9", "D10", "D16", "D19",
    "D20", "D21", "D23", "D27", "D30", "D36", "D37",
    "D39", "D41", "D45", "D46", "D48", "D49", "D51"
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# New cell values, keyed by A1 reference.
$newValues = [ordered]@{
    "D2" = '65.129.34'
    "E2" = '  +0.55%  '
    "D3" = '3.447.05'
    "E3" = '  -1.11%  '
    "E4" = '  +0.01%  '
    "D5" = '576.75'
    "E5" = '  -1.53%  '
    "D6" = '160.78'
    "E6" = '  +0.62%  '
    "D7" = '1.00'
    "E7" = '  +0.06%  '
    "D8" = '3.448.51'
    "E8" = '  -1.21%  '
    "D9" = '0.578'
    "E9" = '  +7.61%  '
    "D10" = '7.26'
    "E10" = '  -4.85%  '
    "E11" = '  -0.09%  '
    "E12" = '  -1.61%  '
    "D13" = '4.043.75'
    "E13" = '  -1.07%  '
    "E15" = '  +0.11%  '
    "D16" = '28.06'
    "E16" = '  +0.99%  '
    "D17" = '65.081.43'
    "E17" = '  +0.47%  '
    "D18" = '3.461.99'
    "E18" = '  -0.78%  '
    "D19" = '6.34'
    "E19" = '  -2.46%  '
    "D20" = '14.23'
    "E20" = '  -1.25%  '
    "D21" = '387.44'
    "E21" = '  -3.21%  '
    "E22" = '  -4.66%  '
    "D23" = '73.09'
    "E23" = '  +0.92%  '
    "E24" = '  -0.98%  '
    "E25" = '  +0.05%  '
    "E26" = '  +9.02%  '
    "D27" = '9.64'
    "E27" = '  -0.71%  '
    "E28" = '  -1.39%  '
    "E29" = '  -0.02%  '
    "D30" = '6.25'
    "E30" = '  +5.66%  '
    "E31" = '  +0.72%  '
    "E32" = '  -1.14%  '
    "E33" = '  -1.97%  '
    "E34" = '  -1.29%  '
    "E35" = '  +1.04%  '
    "B36" = 'Monero'
    "C36" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    "D36" = '161.77'
    "E36" = '  +2.08%  '
    "B37" = 'ImmutableX'
    "C37" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D37" = '1.49'
    "E37" = '  -1.60%  '
    "B38" = 'Maker'
    "C38" = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    "D38" = '3.040.52'
    "E38" = '  +4.25%  '
    "B39" = 'Stacks'
    "C39" = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    "D39" = '1.91'
    "E39" = '  +0.01%  '
    "E40" = '  -2.02%  '
    "D41" = '27.17'
    "E41" = '  -5.21%  '
    "E42" = '  +2.05%  '
    "E43" = '  +1.80%  '
    "E44" = '  -2.07%  '
    "D45" = '0.770'
    "E45" = '  -2.27%  '
    "D46" = '24.88'
    "E46" = '  +7.70%  '
    "E47" = '  -2.81%  '
    "D48" = '2.21'
    "E48" = '  +1.76%  '
    "D49" = '0.868'
    "E49" = '  +2.45%  '
    "E50" = '  +1.57%  '
    "D51" = '305.00'
    "E51" = '  +1.49%  '
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
